$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (Maryland) ---
$ws.Range("B17").Value = Get-Date -Year 2020 -Month 7 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("B17").NumberFormat = "YYYY-MM-DD"

$ws.Range("C17").Value = 74260
$ws.Range("D17").Value = 3202
$ws.Range("E17").Value = 21525
$ws.Range("F17").Value = 1301
$ws.Range("G17").Value = 35.07
$ws.Range("H17").Value = 40.87

$ws.Range("K17").Value = 61384
$ws.Range("L17").Value = 3183

$ws.Range("O17").Value = "Success!"

# --- Row 36 (Iowa) ---
$ws.Range("C36").Value = 35866
